$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row that contains "아몰라그냥해" (row 22),
# shifting all subsequent rows up by one.
$ws.Rows.Item(22).Delete()
